$p = $ppt.ActivePresentation
$new = $p.Slides.Add($p.Slides.Count + 1, 1)
$new.Shapes.Item(1).TextFrame.TextRange.Text = "ICDL"
$new.Shapes.Item(2).TextFrame.TextRange.Text = "Power Point"
